# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
# Rows 3,4,6,8,10,11,12 in column F get incremented to their new values.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1459
    4  = 954
    6  = 2151
    8  = 1309
    10 = 129
    11 = 41
    12 = 315
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
